# Documentation: Added changes from TAFE and screenshots
#
# The "Test Schedule" section of the test plan had five blank Arial/16pt
# paragraphs sitting between the "Testing Tools" list and the results
# table. This fills in four of them with the "Test Objective:" / body,
# "Test Schedule:" / body headings+text (the two headings are italic),
# and removes the stray trailing empty paragraph that used to sit right
# before the table.

$d = $word.ActiveDocument

# Locate the first of the five blank paragraphs by anchoring off the
# last line of the preceding "Testing Tools:" list, rather than a
# hard-coded paragraph index.
$anchor = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "The testing environments labelled above.`r") {
        $anchor = $i
        break
    }
}
$base = $anchor + 1

# Paragraph $base (32): "Test Objective:" heading (italic)
$p = $d.Paragraphs.Item($base)
$r = $p.Range
$r.InsertAfter("Test Objective:")
$r2 = $p.Range
$r2.Font.Name = "Arial"
$r2.Font.NameAscii = "Arial"
$r2.Font.NameBi = "Arial"
$r2.Font.NameOther = "Arial"
$r2.Font.Size = 16
$r2.Font.Italic = $true

# Paragraph $base+1 (33): objective body text (not italic)
$p = $d.Paragraphs.Item($base + 1)
$r = $p.Range
$r.InsertAfter("Discover any erroneous components of code – functionality or experience wise.")
$r2 = $p.Range
$r2.Font.Name = "Arial"
$r2.Font.NameAscii = "Arial"
$r2.Font.NameBi = "Arial"
$r2.Font.NameOther = "Arial"
$r2.Font.Size = 16

# Paragraph $base+2 (34): "Test Schedule:" heading (italic)
$p = $d.Paragraphs.Item($base + 2)
$r = $p.Range
$r.InsertAfter("Test Schedule:")
$r2 = $p.Range
$r2.Font.Name = "Arial"
$r2.Font.NameAscii = "Arial"
$r2.Font.NameBi = "Arial"
$r2.Font.NameOther = "Arial"
$r2.Font.Size = 16
$r2.Font.Italic = $true

# Paragraph $base+3 (35): schedule body text (not italic)
$p = $d.Paragraphs.Item($base + 3)
$r = $p.Range
$r.InsertAfter("After every major feature is added, and extensively before each major release.")
$r2 = $p.Range
$r2.Font.Name = "Arial"
$r2.Font.NameAscii = "Arial"
$r2.Font.NameBi = "Arial"
$r2.Font.NameOther = "Arial"
$r2.Font.Size = 16

# Paragraph $base+4 (36) stays blank. Paragraph $base+5 (37) is the
# stray trailing empty paragraph right before the results table --
# remove it entirely.
$p = $d.Paragraphs.Item($base + 5)
$p.Range.Delete()
